# "test again, no dismantling"
# Flip the two "realistic candidate capacities" switches on the
# "Coupling Parameters" sheet from FALSE to TRUE, which in turn makes the
# dependent IF() formula in C20 recompute to its TRUE-branch text, and
# leave the sheet's active-cell selection on B21 (where the user's cursor
# ended up after making the edits).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Coupling Parameters")
$ws.Activate() | Out-Null

# realistic_candidate_capacities_tobe_installed -> TRUE
$ws.Range("B19").Value = $true

# realistic_candidate_capacities_to_test -> TRUE
# (drives C20's IF(B20=TRUE, ..., ...) formula to the first branch)
$ws.Range("B20").Value = $true

$excel.Calculate() | Out-Null

# Leave the selection where the user left it after editing row 20.
$ws.Range("B21").Select() | Out-Null
